{"js": "// 1. Update the getServices description text: the \"active\" filter option is\n//    being removed, so drop it from the description and tidy the wording.\nconst body = context.document.body;\nconst searchResults = body.search(\n  \"This service, filter services by active, name, max price or returns all services\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"This service, filter services by name and/or max price or returns all services\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 2. Remove the \"active\" row from the getServices parameters table (the\n//    second table in the document: Field/Description rows for name, price,\n//    active).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst paramsTable = tables.items[1];\nparamsTable.rows.load(\"items\");\nawait context.sync();\n\nconst lastRow = paramsTable.rows.items[paramsTable.rows.items.length - 1];\nlastRow.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the getServices description: remove \"active\" as a filter option,\n#    since the \"active\" filter field itself is being removed (see below).\n$find = $d.Content.Find\n$find.Text = \"This service, filter services by active, name, max price or returns all services\"\n$find.Replacement.Text = \"This service, filter services by name and/or max price or returns all services\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2. Remove the \"active\" row from the getServices parameters table (2nd table\n#    in the document: Field/Description for name, price, active).\n$t = $d.Tables.Item(2)\n$lastRow = $t.Rows.Item($t.Rows.Count)\n$lastRow.Delete()\n"}
